$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for rows 2-3
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-15 12:13:10"
$wsOverview.Range("G3").Value = "2016-08-15 12:13:10"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-15 12:13:00"
$wsZhCn.Range("H3").Value = "2016-08-15 12:13:00"
$wsZhCn.Range("K2").Value = "2016-08-15 12:13:27"
$wsZhCn.Range("K3").Value = "2016-08-15 12:13:27"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (K)
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-15 12:13:10"
$wsDeDe.Range("H3").Value = "2016-08-15 12:13:10"
$wsDeDe.Range("K2").Value = "2016-08-15 12:13:34"
$wsDeDe.Range("K3").Value = "2016-08-15 12:13:34"
